$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.528.11'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '2.477.74'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.30'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '92.58'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.545'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.09%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '32.64'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.28%  '
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('E12').Value = '  +2.28%  '
$ws.Range('D13').Value = '2.859.11'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '16.24'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +9.08%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.84'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').Value = '2.476.08'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.768'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.86%  '
$ws.Range('D18').Value = '41.527.21'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.46'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  +2.24%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.98'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +5.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.24'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.61'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.77'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.13%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.62'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.77'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '158.37'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.42'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '17.26'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.00%  '
$ws.Range('B36').Value = 'ApeXProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.43'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -8.71%  '
$ws.Range('E37').Value = '  +4.35%  '
$ws.Range('E38').Value = '  -5.00%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.82'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.77%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -4.09%  '
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').Value = '1.983.34'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '19.32'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.43%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  -2.95%  '
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('D48').Value = '2.717.93'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '97.41'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '67.99'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '72.26'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.50%  '
